$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 35, shifting old rows 35-115 down to 37-117.
$ws.Rows("35:36").Insert()

$oAcute = [char]0xF3
$iAcute = [char]0xED

# New row 35: Early Burlat / Primera record.
$ws.Range("A35").Value = 11
$ws.Range("B35").Value = "Vega Monumental Concepci" + $oAcute + "n"
$ws.Range("C35").Value = "B" + $iAcute + "ob" + $iAcute + "o"
$ws.Range("D35").Value = 44883
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100103
$ws.Range("H35").Value = "Frutos de hueso (carozo)"
$ws.Range("I35").Value = 100103001
$ws.Range("J35").Value = "Cereza"
$ws.Range("K35").Value = "Early Burlat"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 240
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 22000
$ws.Range("P35").Value = 21250
$ws.Range("Q35").Value = "`$/bandeja 10 kilos"
$ws.Range("R35").Value = "Provincia de Curic" + $oAcute
$ws.Range("S35").Value = 2125
$ws.Range("T35").Value = 10

# New row 36: Early Burlat / Segunda record.
$ws.Range("A36").Value = 11
$ws.Range("B36").Value = "Vega Monumental Concepci" + $oAcute + "n"
$ws.Range("C36").Value = "B" + $iAcute + "ob" + $iAcute + "o"
$ws.Range("D36").Value = 44883
$ws.Range("E36").Value = 8
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100103
$ws.Range("H36").Value = "Frutos de hueso (carozo)"
$ws.Range("I36").Value = 100103001
$ws.Range("J36").Value = "Cereza"
$ws.Range("K36").Value = "Early Burlat"
$ws.Range("L36").Value = "Segunda"
$ws.Range("M36").Value = 170
$ws.Range("N36").Value = 12000
$ws.Range("O36").Value = 14000
$ws.Range("P36").Value = 13059
$ws.Range("Q36").Value = "`$/bandeja 10 kilos"
$ws.Range("R36").Value = "Provincia de Curic" + $oAcute
$ws.Range("S36").Value = 1306
$ws.Range("T36").Value = 10

# Ensure date formatting (style) on the new D cells matches the rest of column D.
$ws.Range("D35:D36").NumberFormat = $ws.Range("D34").NumberFormat
